$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.695.76'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '2.043.03'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.16'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.610'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.40'
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("E10").Value = '  +3.28%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '2.348.84'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.34'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.36'
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.48'
$ws.Range("E15").Value = '  +6.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.760'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '2.035.97'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = '37.684.22'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.35'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.89'
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("D21").Value = '0.0₃0827'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.10'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  +2.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.59'
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.73'
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").Value = '  +8.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.36'
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0601'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.45'
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("E37").Value = '  +3.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.46'
$ws.Range("E38").Value = '  +7.13%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.22'
$ws.Range("E40").Value = '  +8.94%  '
$ws.Range("D41").Value = '1.525.35'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.56'
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("E43").Value = '  -1.82%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.19'
$ws.Range("E45").Value = '  +2.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0889'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.94'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.07'
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").Value = '2.237.08'
$ws.Range("E51").Value = '  +0.82%  '
